# Applies the "Jin Chan's Pond of Riches" review edit:
#   - straightens the curly apostrophe / drops " | Review" from the title
#     (both the Heading1 and the later bold "title" paragraph)
#   - reshuffles / edits the "What we like" bullet list
#   - tweaks two "What we don't like" bullets
#   - rewrites the italic meta-description paragraph
#
# NOTE: this runtime's Find.Execute(...,Replace:=) silently "smart-quotes"
# any straight apostrophe in the replacement text, so instead of passing
# a replacement string to Find.Execute we only use Find.Execute to locate
# (and narrow) a Range, then assign Range.Text directly - that path leaves
# straight apostrophes alone.

function Find-TextRange($d, $text) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true)
    if (-not $found) {
        throw "Text not found: $text"
    }
    return $rng
}

function Replace-Text($d, $old, $new) {
    $rng = Find-TextRange $d $old
    $rng.Text = $new
}

function Delete-ParaContaining($d, $text) {
    $rng = Find-TextRange $d $text
    $p = $rng.Paragraphs(1)
    $p.Range.Delete()
}

function Insert-ListItemBefore($d, $anchorText, $newText) {
    $rng = Find-TextRange $d $anchorText
    $p = $rng.Paragraphs(1)
    $idx = $p.Index
    $p.Range.InsertParagraphBefore()
    $newp = $d.Paragraphs($idx)
    $newp.Range.Text = $newText
}

$d = $word.ActiveDocument

# 1. Main heading: curly apostrophe -> straight apostrophe, drop " | Review"
Replace-Text $d "Play Jin Chan’s Pond of Riches for Free | Review" "Play Jin Chan's Pond of Riches for Free"

# 2. "What we like": add a new first bullet about betting options.
Insert-ListItemBefore $d "Impressive graphics and Chinese folklore theme" "Wide range of betting options"

# 3. "What we like": reword the betting-range bullet into the Sticky Respins bullet.
Replace-Text $d "Wide range of bets available" "Sticky Respins feature for additional free spins"

# 4. "What we like": reword the old Sticky Respins bullet into the demo-version bullet.
Replace-Text $d "Sticky Respins feature for extra free spins and prizes" "Availability of a free demo version"

# 5. "What we like": the old standalone demo-version bullet is now redundant - remove it.
Delete-ParaContaining $d "Demo version available to play for free without registration"

# 6. "What we don't like": tighten the wording.
Replace-Text $d "High volatility makes it difficult to win" "High volatility makes winning difficult"

# 7. "What we don't like": drop the explicit percentage.
Replace-Text $d "Low theoretical return to the player (94.27%)" "Low theoretical return to the player"

# 8. Bold "title" paragraph near the end - same fix as step 1.
Replace-Text $d "Play Jin Chan’s Pond of Riches for Free | Review" "Play Jin Chan's Pond of Riches for Free"

# 9. Italic meta-description paragraph.
Replace-Text $d "Read our review of Jin Chan’s Pond of Riches and play this Chinese-themed online slot game for free. Includes pros, cons, and features of the game." "Read our review of Jin Chan's Pond of Riches and play for free without registration."

Write-Output "Edits applied."
